# Generate Report for Archive
# - Flip status text "Ready for handoff" -> "In Translation" everywhere it
#   appears (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# - Shrink the now-narrower "Status" columns to match the new text's width
#   (Overview columns E:F, zh-cn column C, de-de column C).

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "In Translation"

# Closest achievable ColumnWidth input for the target stored width of
# ~13.41 characters (this host snaps stored widths to an MDW-6 pixel grid,
# i.e. stored = (Round(input*6)+5)/6, so 12.5 -> 13.333333333333334, the
# nearest representable width to the target).
$newColWidth = 12.5

# --- Overview sheet: columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Columns("E:F").ColumnWidth = $newColWidth

# --- zh-cn sheet: Status column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Columns("C").ColumnWidth = $newColWidth

# --- de-de sheet: Status column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Columns("C").ColumnWidth = $newColWidth
